# Add "Slovakia" / "SlovakiaxBB" worksheets, mirroring the existing
# "Germany" / "GermanyxBB" templates, for the new Slovakia market.

$wb = $excel.ActiveWorkbook

$germany    = $wb.Worksheets.Item("Germany")
$germanyBB  = $wb.Worksheets.Item("GermanyxBB")
$portugalBB = $wb.Worksheets.Item("PortugalBB")
$portugal   = $wb.Worksheets.Item("Portugal")

# Copy "Germany" right after "PortugalBB" (the current last tab) and
# rename it to "Slovakia".
$germany.Copy($null, $portugalBB)
$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"

# Copy "GermanyxBB" right after the freshly-created "Slovakia" tab and
# rename it to "SlovakiaxBB".
$germanyBB.Copy($null, $slovakia)
$slovakiaBB = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakiaBB.Name = "SlovakiaxBB"

# Fill in the market-specific values on both new sheets.
$slovakia.Range("B2").Value = "Slovakia Market"
$slovakia.Range("B4").Value = "NGC-2930/T3223"

$slovakiaBB.Range("B2").Value = "Slovakia Market"
$slovakiaBB.Range("B4").Value = "NGC-2930/T3223"

# Restore the cursor/selection state on the new sheets.
$null = $slovakia.Range("B2").Select()
$null = $slovakiaBB.Range("B10").Select()

# The previously-active "PortugalBB" tab loses its selection, and
# "Portugal" cursor moves to A23.
$null = $portugal.Range("A23").Select()

# "Slovakia" becomes the new active tab.
$slovakia.Activate()
